# Updated cryptos list on Sat Jul 27 14:26:10 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.051.95'
$ws.Range('E2').Value = '  +2.16%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.307.38'
$ws.Range('E3').Value = '  +1.63%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.80'
$ws.Range('E5').Value = '  +2.07%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.06'
$ws.Range('E6').Value = '  +2.44%  '

# Row 7
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').Value = '  +1.92%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.136'
$ws.Range('E9').Value = '  +4.53%  '

# Row 10
$ws.Range('E10').Value = '  -1.03%  '

# Row 11
$ws.Range('E11').Value = '  +2.67%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.885.52'
$ws.Range('E12').Value = '  +2.21%  '

# Row 13
$ws.Range('E13').Value = '  +0.27%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.33'
$ws.Range('E14').Value = '  +4.94%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.074.28'
$ws.Range('E15').Value = '  +2.21%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('E16').Value = '  +3.64%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.284.54'
$ws.Range('E17').Value = '  +0.85%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.89'
$ws.Range('E18').Value = '  +1.21%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.77'
$ws.Range('E19').Value = '  +2.57%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '389.45'
$ws.Range('E20').Value = '  +4.25%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.79'
$ws.Range('E21').Value = '  +2.65%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.06'
$ws.Range('E22').Value = '  +1.44%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.28%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000123'
$ws.Range('E24').Value = '  +3.23%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.520'
$ws.Range('E25').Value = '  +1.94%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.80'
$ws.Range('E26').Value = '  +2.09%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.188'
$ws.Range('E27').Value = '  +4.20%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.27%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.84'
$ws.Range('E29').Value = '  +3.18%  '

# Row 30
$ws.Range('E30').Value = '  +1.66%  '

# Row 31
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.14'
$ws.Range('E31').Value = '  +2.48%  '

# Row 32
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.32'
$ws.Range('E32').Value = '  +3.86%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.19'
$ws.Range('E33').Value = '  +4.72%  '

# Row 34
$ws.Range('E34').Value = '  +0.00%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  +4.14%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.40'
$ws.Range('E36').Value = '  +0.81%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.92'
$ws.Range('E37').Value = '  +3.52%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.842'
$ws.Range('E38').Value = '  -1.97%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.80'
$ws.Range('E39').Value = '  +0.01%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.63'
$ws.Range('E40').Value = '  +5.38%  '

# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.71'
$ws.Range('E41').Value = '  -1.20%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.64'
$ws.Range('E42').Value = '  +1.83%  '

# Row 43
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.73'
$ws.Range('E43').Value = '  +2.87%  '

# Row 44
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0697'
$ws.Range('E44').Value = '  +3.47%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.55'
$ws.Range('E45').Value = '  -0.23%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.652.44'
$ws.Range('E46').Value = '  -2.03%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '339.99'
$ws.Range('E47').Value = '  -5.85%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0286'
$ws.Range('E48').Value = '  +2.81%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.48'
$ws.Range('E49').Value = '  +5.42%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.28%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.33'
$ws.Range('E51').Value = '  +3.63%  '

